# Update profit-calculation figures on the Leve profit sheets (scheduled market-price refresh).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 158.4
$ws.Range("I12").Value = 158.4
$ws.Range("K12").Value = 158.4
$ws.Range("M12").Value = 11.59999999999999
$ws.Range("H19").Value = 718.5
$ws.Range("I19").Value = 693.5
$ws.Range("J19").Value = 735.1667
$ws.Range("K19").Value = 693.5
$ws.Range("L19").Value = 735.1667
$ws.Range("M19").Value = -518.5
$ws.Range("N19").Value = -1085.1667
$ws.Range("H33").Value = 189.23529
$ws.Range("I33").Value = 112.07143
$ws.Range("J33").Value = 549.3333
$ws.Range("K33").Value = 112.07143
$ws.Range("L33").Value = 549.3333
$ws.Range("M33").Value = 116.92857
$ws.Range("N33").Value = -1007.3333
$ws.Range("H62").Value = 4255.4443
$ws.Range("I62").Value = 3537.375
$ws.Range("K62").Value = 3537.375
$ws.Range("M62").Value = -2913.375
$ws.Range("H65").Value = 4255.4443
$ws.Range("I65").Value = 3537.375
$ws.Range("K65").Value = 17686.875
$ws.Range("M65").Value = -14566.875
$ws.Range("H107").Value = 9685.714
$ws.Range("I107").Value = 9600
$ws.Range("J107").Value = 9900
$ws.Range("K107").Value = 9600
$ws.Range("L107").Value = 9900
$ws.Range("M107").Value = -7680
$ws.Range("N107").Value = -13740
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 32702
$ws.Range("J110").Value = 32702
$ws.Range("L110").Value = 32702
$ws.Range("N110").Value = -40882
$ws.Range("H111").Value = 3550
$ws.Range("I111").Value = 3660
$ws.Range("K111").Value = 10980
$ws.Range("M111").Value = -7913
$ws.Range("H113").Value = 2557.6
$ws.Range("I113").Value = 2557.6
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2557.6
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 696.4000000000001
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 5470127.5
$ws.Range("I132").Value = 7411348.5
$ws.Range("K132").Value = 22234045.5
$ws.Range("M132").Value = -22231515.5
$ws.Range("H137").Value = 1056.2593
$ws.Range("I137").Value = 886.875
$ws.Range("J137").Value = 1221.5122
$ws.Range("K137").Value = 2660.625
$ws.Range("L137").Value = 3664.536599999999
$ws.Range("M137").Value = -110.625
$ws.Range("N137").Value = -8764.536599999999
$ws.Range("H138").Value = 1170.02
$ws.Range("J138").Value = 1701.875
$ws.Range("L138").Value = 5105.625
$ws.Range("N138").Value = -15385.625
$ws.Range("H141").Value = 613.4727
$ws.Range("I141").Value = 512.2308
$ws.Range("K141").Value = 1536.6924
$ws.Range("M141").Value = 3643.3076

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3678.7693
$ws.Range("I32").Value = 3488.2258
$ws.Range("J32").Value = 7616.6665
$ws.Range("K32").Value = 3488.2258
$ws.Range("L32").Value = 7616.6665
$ws.Range("M32").Value = -3201.2258
$ws.Range("N32").Value = -8190.6665
$ws.Range("H45").Value = 1566.6364
$ws.Range("I45").Value = 1543.3
$ws.Range("K45").Value = 1543.3
$ws.Range("M45").Value = -1166.3
$ws.Range("H112").Value = 7846.75
$ws.Range("J112").Value = 7846.75
$ws.Range("L112").Value = 7846.75
$ws.Range("N112").Value = -10800.75
$ws.Range("H122").Value = 3499.8
$ws.Range("I122").Value = 3874.75
$ws.Range("K122").Value = 11624.25
$ws.Range("M122").Value = -9174.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 27780330
$ws.Range("I94").Value = 31252496
$ws.Range("K94").Value = 31252496
$ws.Range("M94").Value = -31252045
$ws.Range("H99").Value = 333335330
$ws.Range("I99").Value = 500001500
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 500001500
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -500000002
$ws.Range("N99").Value = -5996
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22372
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -71856
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 1080.76
$ws.Range("I134").Value = 1147.8
$ws.Range("J134").Value = 812.6
$ws.Range("K134").Value = 3443.4
$ws.Range("L134").Value = 2437.8
$ws.Range("M134").Value = -908.3999999999996
$ws.Range("N134").Value = -7507.8
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2000.2222
$ws.Range("I5").Value = 2221.7144
$ws.Range("J5").Value = 1225
$ws.Range("K5").Value = 6665.1432
$ws.Range("L5").Value = 3675
$ws.Range("M5").Value = -6553.1432
$ws.Range("N5").Value = -3899
$ws.Range("H68").Value = 1310.3334
$ws.Range("I68").Value = 1286.5454
$ws.Range("K68").Value = 3859.6362
$ws.Range("M68").Value = -3048.6362
$ws.Range("H71").Value = 1310.3334
$ws.Range("I71").Value = 1286.5454
$ws.Range("K71").Value = 11578.9086
$ws.Range("M71").Value = -7522.908599999999
$ws.Range("H81").Value = 2559.0715
$ws.Range("I81").Value = 1503
$ws.Range("K81").Value = 4509
$ws.Range("M81").Value = -3386
$ws.Range("H84").Value = 2559.0715
$ws.Range("I84").Value = 1503
$ws.Range("K84").Value = 13527
$ws.Range("M84").Value = -7911
$ws.Range("H113").Value = 681.1515000000001
$ws.Range("J113").Value = 686.8125
$ws.Range("L113").Value = 2060.4375
$ws.Range("N113").Value = -6400.4375
$ws.Range("H135").Value = 2000.2222
$ws.Range("I135").Value = 2221.7144
$ws.Range("J135").Value = 1225
$ws.Range("K135").Value = 19995.4296
$ws.Range("L135").Value = 11025
$ws.Range("M135").Value = -17460.4296
$ws.Range("N135").Value = -16095

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1763.25
$ws.Range("I7").Value = 1625.4
$ws.Range("J7").Value = 2452.5
$ws.Range("K7").Value = 1625.4
$ws.Range("L7").Value = 2452.5
$ws.Range("M7").Value = -1513.4
$ws.Range("N7").Value = -2676.5
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 20843484
$ws.Range("I122").Value = 41684630
$ws.Range("J122").Value = 2333.1667
$ws.Range("K122").Value = 125053890
$ws.Range("L122").Value = 6999.500100000001
$ws.Range("M122").Value = -125051440
$ws.Range("N122").Value = -11899.5001
$ws.Range("H126").Value = 1763.25
$ws.Range("I126").Value = 1625.4
$ws.Range("J126").Value = 2452.5
$ws.Range("K126").Value = 4876.200000000001
$ws.Range("L126").Value = 7357.5
$ws.Range("M126").Value = -2406.200000000001
$ws.Range("N126").Value = -12297.5
$ws.Range("H132").Value = 22833.334
$ws.Range("I132").Value = 1553.5385
$ws.Range("J132").Value = 47982.184
$ws.Range("K132").Value = 4660.6155
$ws.Range("L132").Value = 143946.552
$ws.Range("M132").Value = -2130.6155
$ws.Range("N132").Value = -149006.552

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 178574660
$ws.Range("I122").Value = 250002510
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 750007530
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -750005080
$ws.Range("N122").Value = -19898.5
$ws.Range("H132").Value = 2689.2368
$ws.Range("I132").Value = 2645.394
$ws.Range("J132").Value = 2978.6
$ws.Range("K132").Value = 7936.181999999999
$ws.Range("L132").Value = 8935.799999999999
$ws.Range("M132").Value = -5406.181999999999
$ws.Range("N132").Value = -13995.8
$ws.Range("H136").Value = 540.81396
$ws.Range("I136").Value = 380.51724
$ws.Range("J136").Value = 872.8570999999999
$ws.Range("K136").Value = 1141.55172
$ws.Range("L136").Value = 2618.5713
$ws.Range("M136").Value = 1408.44828
$ws.Range("N136").Value = -7718.5713
